$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 154.18182
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# Row 4
$ws.Range("H4").Value = 208.5
$ws.Range("I4").Value = 208.5
$ws.Range("K4").Value = 208.5
$ws.Range("M4").Value = -94.5
# Row 18
$ws.Range("H18").Value = 950
$ws.Range("I18").Value = 950
$ws.Range("K18").Value = 950
$ws.Range("M18").Value = -666
# Row 33
$ws.Range("H33").Value = 392.82608
$ws.Range("I33").Value = 113.85714
$ws.Range("K33").Value = 113.85714
$ws.Range("M33").Value = 115.14286
# Row 38
$ws.Range("H38").Value = 347.33334
$ws.Range("I38").Value = 106.181816
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 318.545448
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 53.45455200000004
$ws.Range("N38").Value = -9744
# Row 58
$ws.Range("H58").Value = 400
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 82
$ws.Range("H82").Value = 1282
$ws.Range("I82").Value = 1282
$ws.Range("K82").Value = 3846
$ws.Range("M82").Value = -3440
# Row 85
$ws.Range("H85").Value = 1282
$ws.Range("I85").Value = 1282
$ws.Range("K85").Value = 3846
$ws.Range("M85").Value = -2442
# Row 87
$ws.Range("H87").Value = 64995
$ws.Range("J87").Value = 64995
$ws.Range("L87").Value = 64995
$ws.Range("N87").Value = -67491
# Row 90
$ws.Range("H90").Value = 64995
$ws.Range("J90").Value = 64995
$ws.Range("L90").Value = 194985
$ws.Range("N90").Value = -207465
# Row 135
$ws.Range("H135").Value = 1532.8572
$ws.Range("I135").Value = 1660
$ws.Range("K135").Value = 14940
$ws.Range("M135").Value = -12405

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 5928
$ws.Range("J2").Value = 6506.5
$ws.Range("L2").Value = 6506.5
$ws.Range("N2").Value = -6732.5
# Row 61
$ws.Range("H61").Value = 7097981.5
$ws.Range("I61").Value = 8337922
$ws.Range("K61").Value = 8337922
$ws.Range("M61").Value = -8337710
# Row 74
$ws.Range("H74").Value = 4060.1538
$ws.Range("J74").Value = 7799.5713
$ws.Range("L74").Value = 7799.5713
$ws.Range("N74").Value = -9547.5713
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 77
$ws.Range("H77").Value = 4060.1538
$ws.Range("J77").Value = 7799.5713
$ws.Range("L77").Value = 38997.85649999999
$ws.Range("N77").Value = -47733.85649999999
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 94
$ws.Range("H94").Value = 12000
$ws.Range("J94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("N94").Value = -13802
# Row 116
$ws.Range("H116").Value = 5928
$ws.Range("J116").Value = 6506.5
$ws.Range("L116").Value = 6506.5
$ws.Range("N116").Value = -11094.5
# Row 136
$ws.Range("H136").Value = 7097981.5
$ws.Range("I136").Value = 8337922
$ws.Range("K136").Value = 25013766
$ws.Range("M136").Value = -25011216

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 5928
$ws.Range("J3").Value = 6506.5
$ws.Range("L3").Value = 6506.5
$ws.Range("N3").Value = -6734.5
# Row 80
$ws.Range("H80").Value = 582.41174
$ws.Range("I80").Value = 269.42856
$ws.Range("J80").Value = 801.5
$ws.Range("K80").Value = 269.42856
$ws.Range("L80").Value = 801.5
$ws.Range("M80").Value = 728.5714399999999
$ws.Range("N80").Value = -2797.5
# Row 83
$ws.Range("H83").Value = 582.41174
$ws.Range("I83").Value = 269.42856
$ws.Range("J83").Value = 801.5
$ws.Range("K83").Value = 1347.1428
$ws.Range("L83").Value = 4007.5
$ws.Range("M83").Value = 3644.8572
$ws.Range("N83").Value = -13991.5
# Row 134
$ws.Range("H134").Value = 7359.857
$ws.Range("I134").Value = 7378.8184
$ws.Range("K134").Value = 22136.4552
$ws.Range("M134").Value = -19601.4552

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 477.14285
$ws.Range("I7").Value = 188.18182
$ws.Range("J7").Value = 1536.6666
$ws.Range("K7").Value = 188.18182
$ws.Range("L7").Value = 1536.6666
$ws.Range("M7").Value = -75.18181999999999
$ws.Range("N7").Value = -1762.6666
# Row 22
$ws.Range("H22").Value = 7751.5
$ws.Range("I22").Value = 7751.5
$ws.Range("K22").Value = 7751.5
$ws.Range("M22").Value = -7401.5
# Row 31
$ws.Range("H31").Value = 34488236
$ws.Range("I31").Value = 62504572
$ws.Range("K31").Value = 62504572
$ws.Range("M31").Value = -62504277
# Row 34
$ws.Range("H34").Value = 34488236
$ws.Range("I34").Value = 62504572
$ws.Range("K34").Value = 62504572
$ws.Range("M34").Value = -62504370
# Row 122
$ws.Range("H122").Value = 5436.273
$ws.Range("I122").Value = 5255.4443
$ws.Range("K122").Value = 15766.3329
$ws.Range("M122").Value = -13316.3329
# Row 132
$ws.Range("H132").Value = 3385.7144
$ws.Range("I132").Value = 4020
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 12060
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -9530
$ws.Range("N132").Value = -10460
# Row 134
$ws.Range("I134").Value = 2696.077
$ws.Range("J134").Value = 10366
$ws.Range("K134").Value = 8088.231000000001
$ws.Range("L134").Value = 31098
$ws.Range("M134").Value = -5553.231000000001
$ws.Range("N134").Value = -36168

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 5440.263
$ws.Range("J2").Value = 11305.556
$ws.Range("L2").Value = 67833.33600000001
$ws.Range("N2").Value = -68059.33600000001
# Row 5
$ws.Range("H5").Value = 3656.5
$ws.Range("I5").Value = 3007.8
$ws.Range("J5").Value = 6900
$ws.Range("K5").Value = 9023.400000000001
$ws.Range("L5").Value = 20700
$ws.Range("M5").Value = -8911.400000000001
$ws.Range("N5").Value = -20924
# Row 12
$ws.Range("H12").Value = 233.71428
$ws.Range("I12").Value = 20.625
$ws.Range("J12").Value = 517.8333
$ws.Range("K12").Value = 61.875
$ws.Range("L12").Value = 1553.4999
$ws.Range("M12").Value = 111.125
$ws.Range("N12").Value = -1899.4999
# Row 37
$ws.Range("H37").Value = 61786.156
$ws.Range("J37").Value = 61786.156
$ws.Range("L37").Value = 185358.468
$ws.Range("N37").Value = -185582.468
# Row 70
$ws.Range("H70").Value = 400
$ws.Range("I70").Value = 400
$ws.Range("K70").Value = 1200
$ws.Range("M70").Value = -885
# Row 73
$ws.Range("H73").Value = 400
$ws.Range("I73").Value = 400
$ws.Range("K73").Value = 1200
$ws.Range("M73").Value = -108
# Row 109
$ws.Range("H109").Value = 5624.875
$ws.Range("I109").Value = 2499
$ws.Range("J109").Value = 6071.4287
$ws.Range("K109").Value = 7497
$ws.Range("L109").Value = 18214.2861
$ws.Range("M109").Value = -6457
$ws.Range("N109").Value = -20294.2861
# Row 125
$ws.Range("H125").Value = 17250
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840
# Row 129
$ws.Range("H129").Value = 6235.4707
$ws.Range("J129").Value = 7916.9
$ws.Range("L129").Value = 23750.7
$ws.Range("N129").Value = -33750.7
# Row 135
$ws.Range("H135").Value = 3656.5
$ws.Range("I135").Value = 3007.8
$ws.Range("J135").Value = 6900
$ws.Range("K135").Value = 27070.2
$ws.Range("L135").Value = 62100
$ws.Range("M135").Value = -24535.2
$ws.Range("N135").Value = -67170

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2630.7693
$ws.Range("I126").Value = 1525
$ws.Range("K126").Value = 4575
$ws.Range("M126").Value = -2105
# Row 132
$ws.Range("H132").Value = 2194.0571
$ws.Range("I132").Value = 1792.7858
$ws.Range("K132").Value = 5378.357400000001
$ws.Range("M132").Value = -2848.357400000001
# Row 140
$ws.Range("H140").Value = 10000
$ws.Range("J140").Value = 10000
$ws.Range("L140").Value = 10000
$ws.Range("N140").Value = -20360

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 140
$ws.Range("H140").Value = 94623
$ws.Range("J140").Value = 94623
$ws.Range("L140").Value = 94623
$ws.Range("N140").Value = -104983

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5752.12
$ws.Range("I132").Value = 5121.1055
$ws.Range("J132").Value = 7750.3335
$ws.Range("K132").Value = 15363.3165
$ws.Range("L132").Value = 23251.0005
$ws.Range("M132").Value = -12833.3165
$ws.Range("N132").Value = -28311.0005
